# Updated cryptos list values (Price column D, Volume(1h) column E)
# generated from the target OOXML diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.943.35'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.332.00'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +4.39%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '97.53'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '271.86'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.98%  '
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.627'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '46.09'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('E11').Value = '  +3.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.14'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.682.60'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.66'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.45%  '
$ws.Range('E16').Value = '  +8.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.334.84'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +4.89%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.895.35'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.76%  '
$ws.Range('E19').Value = '  +5.45%  '
$ws.Range('E20').Value = '  +7.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.92'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '240.33'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.01%  '
$ws.Range('E23').Value = '  -1.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.51'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +5.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.998'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.45'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.76%  '
$ws.Range('E27').Value = '  +1.14%  '
$ws.Range('E28').Value = '  -1.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.29'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.32'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -5.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.48'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +7.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '173.56'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0907'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.55%  '
$ws.Range('E34').Value = '  +0.68%  '
$ws.Range('E35').Value = '  +2.64%  '
$ws.Range('E36').Value = '  +3.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.110'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.51%  '
$ws.Range('E38').Value = '  +2.40%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.40'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -5.68%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.37'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +8.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.239'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +9.61%  '
$ws.Range('E42').Value = '  +18.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.28'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.22'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +10.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.54'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.41'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('E47').Value = '  +5.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '100.55'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.17%  '
$ws.Range('E49').Value = '  +1.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.559.52'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +4.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.188'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +15.73%  '
